$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.64973783493042
$ws.Range("B1").Value = 1.928603172302246
$ws.Range("C1").Value = 2.480955123901367
$ws.Range("D1").Value = 4.155784606933594
$ws.Range("E1").Value = 1.522474050521851
